# Task 27 update
#
# 1. On the "IndexPage" sheet, insert a new column C ("VerifyUser") between
#    the existing "Flow" and "NewSubmission" columns, with a sample
#    "AreEqual>Thomas Charles" value under it for the first test case. This
#    shifts "NewSubmission" to column D and "ServicePolicy" to column E.
# 2. Update the active sheet/selection bookkeeping: IndexPage becomes the
#    active sheet (previously NewSubmissionPage), with C2 selected there,
#    while NewSubmissionPage keeps a plain B4 selection (no longer scrolled
#    to K1 / O12).

$wb = $excel.ActiveWorkbook

$wsIndex  = $wb.Worksheets.Item("IndexPage")
$wsNewSub = $wb.Worksheets.Item("NewSubmissionPage")

# --- 1. Insert the new "VerifyUser" column on IndexPage -------------------
$wsIndex.Columns("C").Insert()

$wsIndex.Range("C1").Value = "VerifyUser"
$wsIndex.Range("C2").Value = "AreEqual>Thomas Charles"

$wsIndex.Columns("C").AutoFit()

# --- 2. Update selections / active sheet -----------------------------------
# Record the new selection on NewSubmissionPage first so that the later
# activation of IndexPage is what "sticks" as the workbook's active tab.
$wsNewSub.Range("B4").Select()

$wsIndex.Activate()
$wsIndex.Range("C2").Select()
